$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (text, qty, price) replacing the old 10-row table with 8 rows.
$data = @(
    @("Работа (от 1000р клад)", 5000, 193),
    @("Tv5 C4", 3200, 350),
    @("Tv10 C4(Готовый)", 6000, 259),
    @("Ck1 White Diamond", 2500, 235),
    @("Ck05 White Diamond", 1500, 261),
    @("Ck5 White Diamond", 10000, 435),
    @("Tv2 C4", 1800, 211),
    @("Ck10 White Diamond", 18000, 434)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# The table shrank from 10 rows to 8 - clear out the now-unused trailing rows.
$ws.Range("A9:C10").Delete()

$ws.Range("C8").Select()
